$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.638.70"
$ws.Range("E2").Value = "  +2.01%  "

# Row 3
$ws.Range("D3").Value = "3.919.74"
$ws.Range("E3").Value = "  +1.37%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").Value = "'480.47"
$ws.Range("E5").Value = "  +2.41%  "

# Row 6
$ws.Range("D6").Value = "'144.29"
$ws.Range("E6").Value = "  -0.46%  "

# Row 7
$ws.Range("D7").Value = "'0.620"
$ws.Range("E7").Value = "  -2.61%  "

# Row 8
$ws.Range("E8").Value = "  -0.15%  "

# Row 9
$ws.Range("E9").Value = "  -3.14%  "

# Row 10
$ws.Range("E10").Value = "  +8.65%  "

# Row 11
$ws.Range("D11").Value = "'0.0000350"
$ws.Range("E11").Value = "  +11.71%  "

# Row 12
$ws.Range("D12").Value = "'42.51"
$ws.Range("E12").Value = "  -2.53%  "

# Row 13
$ws.Range("E13").Value = "  +0.32%  "

# Row 14
$ws.Range("D14").Value = "4.551.38"
$ws.Range("E14").Value = "  +1.12%  "

# Row 15
$ws.Range("B15").Value = "Uniswap"
$ws.Range("C15").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D15").Value = "'14.58"
$ws.Range("E15").Value = "  -1.86%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.911.09"
$ws.Range("E16").Value = "  -0.35%  "

# Row 17
$ws.Range("E17").Value = "  -0.38%  "

# Row 18
$ws.Range("D18").Value = "'19.62"
$ws.Range("E18").Value = "  -2.39%  "

# Row 19
$ws.Range("E19").Value = "  -3.09%  "

# Row 20
$ws.Range("D20").Value = "68.734.48"
$ws.Range("E20").Value = "  +1.94%  "

# Row 21
$ws.Range("D21").Value = "'432.34"
$ws.Range("E21").Value = "  -0.77%  "

# Row 22
$ws.Range("D22").Value = "'14.61"
$ws.Range("E22").Value = "  -2.08%  "

# Row 23
$ws.Range("E23").Value = "  +1.79%  "

# Row 24
$ws.Range("D24").Value = "'87.27"
$ws.Range("E24").Value = "  -2.07%  "

# Row 25
$ws.Range("D25").Value = "'11.61"
$ws.Range("E25").Value = "  +16.26%  "

# Row 26
$ws.Range("D26").Value = "'3.56"
$ws.Range("E26").Value = "  -1.31%  "

# Row 27
$ws.Range("D27").Value = "'38.03"
$ws.Range("E27").Value = "  +0.32%  "

# Row 28
$ws.Range("D28").Value = "'10.16"
$ws.Range("E28").Value = "  -0.38%  "

# Row 29
$ws.Range("E29").Value = "  +6.12%  "

# Row 30
$ws.Range("D30").Value = "'704.63"
$ws.Range("E30").Value = "  -3.39%  "

# Row 31
$ws.Range("E31").Value = "  -3.41%  "

# Row 32
$ws.Range("D32").Value = "'13.23"
$ws.Range("E32").Value = "  -4.53%  "

# Row 33
$ws.Range("D33").Value = "'2.86"
$ws.Range("E33").Value = "  +3.28%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0902"
$ws.Range("E34").Value = "  +30.22%  "

# Row 35
$ws.Range("D35").Value = "'41.29"
$ws.Range("E35").Value = "  -7.61%  "

# Row 36
$ws.Range("D36").Value = "'58.75"
$ws.Range("E36").Value = "  +1.31%  "

# Row 37
$ws.Range("D37").Value = "'0.152"
$ws.Range("E37").Value = "  -6.96%  "

# Row 38
$ws.Range("D38").Value = "'5.62"
$ws.Range("E38").Value = "  +1.21%  "

# Row 39
$ws.Range("E39").Value = "  -0.11%  "

# Row 40
$ws.Range("D40").Value = "'0.0472"
$ws.Range("E40").Value = "  -2.54%  "

# Row 41
$ws.Range("E41").Value = "  +9.48%  "

# Row 42
$ws.Range("D42").Value = "'2.73"
$ws.Range("E42").Value = "  +5.32%  "

# Row 43
$ws.Range("E43").Value = "  +1.97%  "

# Row 44
$ws.Range("E44").Value = "  -3.15%  "

# Row 45
$ws.Range("E45").Value = "  -0.35%  "

# Row 46
$ws.Range("E46").Value = "  +0.01%  "

# Row 47
$ws.Range("E47").Value = "  -1.41%  "

# Row 48
$ws.Range("E48").Value = "  -1.08%  "

# Row 49
$ws.Range("D49").Value = "'147.40"
$ws.Range("E49").Value = "  +2.16%  "

# Row 50
$ws.Range("D50").Value = "'3.15"
$ws.Range("E50").Value = "  -4.07%  "

# Row 51
$ws.Range("E51").Value = "  -2.52%  "
